$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 842, shifting all rows from 842 downward
# (old row 842 "2026/12/29" etc. becomes row 843, ..., old row 883 becomes row 884)
$ws.Rows.Item(842).Insert()

# Fill the new row 842 with the inserted data point.
# Force the date column to be treated as plain text (not auto-parsed into a date
# serial number) by setting the number format to text before assigning the value,
# then restore the cell to the default "Normal" style so no stray formatting
# remains on the cell.
$ws.Range("A842").NumberFormat = "@"
$ws.Range("A842").Value = "2026/02/21"
$ws.Range("A842").Style = "Normal"

$ws.Range("B842").Value = "土"
$ws.Range("C842").Value = 12
$ws.Range("D842").Value = 201
